# Applies the crypto price/volume refresh described in the commit
# "Updated cryptos list on Fri Aug  9 11:53:00 UTC 2024 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.702.46"
$ws.Range("E2").Value = "  +6.15%  "
$ws.Range("D3").Value = "2.637.77"
$ws.Range("E3").Value = "  +8.60%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.996"
$ws.Range("E4").Value = "  -0.60%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "507.80"
$ws.Range("E5").Value = "  +3.74%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "156.69"
$ws.Range("E6").Value = "  +1.44%  "
$ws.Range("E7").Value = "  -0.12%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.587"
$ws.Range("E8").Value = "  -4.82%  "
$ws.Range("D9").Value = "2.630.10"
$ws.Range("E9").Value = "  +7.48%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.41"
$ws.Range("E10").Value = "  +3.56%  "
$ws.Range("E11").Value = "  +4.32%  "
$ws.Range("E12").Value = "  +2.63%  "
$ws.Range("E13").Value = "  +0.99%  "
$ws.Range("D14").Value = "3.050.36"
$ws.Range("E14").Value = "  +6.85%  "
$ws.Range("D15").Value = "60.630.98"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "21.72"
$ws.Range("E16").Value = "  +5.41%  "
$ws.Range("E17").Value = "  +5.05%  "
$ws.Range("D18").Value = "2.622.91"
$ws.Range("E18").Value = "  +7.33%  "
$ws.Range("E19").Value = "  +3.23%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "344.39"
$ws.Range("E20").Value = "  +6.10%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.42"
$ws.Range("E21").Value = "  +3.99%  "
$ws.Range("E22").Value = "  +4.04%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.998"
$ws.Range("E23").Value = "  -0.04%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.74"
$ws.Range("E24").Value = "  -0.29%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "60.24"
$ws.Range("E25").Value = "  +4.07%  "
$ws.Range("E26").Value = "  +5.53%  "
$ws.Range("E27").Value = "  +3.54%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.991"
$ws.Range("E28").Value = "  -0.58%  "
$ws.Range("D29").Value = "0.0₃0855"
$ws.Range("E29").Value = "  +8.69%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.55"
$ws.Range("E30").Value = "  +3.86%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "156.24"
$ws.Range("E32").Value = "  +3.94%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "19.37"
$ws.Range("E33").Value = "  +3.54%  "
$ws.Range("E34").Value = "  +3.22%  "
$ws.Range("E35").Value = "  +7.50%  "
$ws.Range("E36").Value = "  +5.93%  "
$ws.Range("E37").Value = "  +5.98%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "309.83"
$ws.Range("E38").Value = "  +8.40%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.852"
$ws.Range("E39").Value = "  +4.05%  "
$ws.Range("E40").Value = "  +7.35%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.48"
$ws.Range("E41").Value = "  +7.29%  "
$ws.Range("E42").Value = "  +27.38%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "35.68"
$ws.Range("E43").Value = "  +4.94%  "
$ws.Range("B44").Value = "Mantle"
$ws.Range("C44").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.627"
$ws.Range("E44").Value = "  +3.81%  "
$ws.Range("B45").Value = "Hedera"
$ws.Range("C45").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0569"
$ws.Range("E45").Value = "  +6.90%  "
$ws.Range("E46").Value = "  -0.58%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.998"
$ws.Range("E47").Value = "  +0.52%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "19.87"
$ws.Range("E48").Value = "  +12.85%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "4.92"
$ws.Range("E49").Value = "  +7.06%  "
$ws.Range("E50").Value = "  +3.68%  "
$ws.Range("D51").Value = "2.044.70"
$ws.Range("E51").Value = "  +7.91%  "
